$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated vm_pu results for Case_4_44 (380 kV slack voltage changed from 1.05 to 1.02 p.u.)
# New values for rows 2-25, columns B,C,D,E,F,I,J,K,L,M,N (column G stays 1, H is empty)
$data = @{
    2 = @{ "B"=1.02; "C"=1.031084866867711; "D"=1.034459721883535; "E"=1.034723110283632; "F"=1.029680412142678; "I"=1.036335185846188; "J"=1.036222541304671; "K"=1.037259133915187; "L"=1.037521766041777; "M"=1.032493620167627; "N"=1.037694096169729 }
    3 = @{ "B"=1.02; "C"=1.0320807527608; "D"=1.035199337676789; "E"=1.035666960242666; "F"=1.031308047458651; "I"=1.036615465583436; "J"=1.036859917068646; "K"=1.037808127357966; "L"=1.038274504481072; "M"=1.033927248383229; "N"=1.038332377080379 }
    4 = @{ "B"=1.02; "C"=1.032724876323362; "D"=1.035677608407569; "E"=1.03627776384439; "F"=1.032360902234008; "I"=1.036795293688005; "J"=1.037271478711103; "K"=1.038162388674269; "L"=1.038761026308899; "M"=1.034854104806733; "N"=1.038744523187574 }
    5 = @{ "B"=1.02; "C"=1.032995599411111; "D"=1.035878599183474; "E"=1.036534562758361; "F"=1.032803448790433; "I"=1.036870526949768; "J"=1.037444293157647; "K"=1.038311087213265; "L"=1.038965428854925; "M"=1.035243569995311; "N"=1.038917583050455 }
    6 = @{ "B"=1.02; "C"=1.033041051136637; "D"=1.035912342065855; "E"=1.036577681439803; "F"=1.032877750174697; "I"=1.036883137461072; "J"=1.03747329739586; "K"=1.038336040702296; "L"=1.038999741261018; "M"=1.035308952211278; "N"=1.038946628478012 }
    7 = @{ "B"=1.02; "C"=1.032728493998575; "D"=1.03568029434787; "E"=1.036281195136417; "F"=1.032366815843973; "I"=1.036796300398277; "J"=1.037273788675217; "K"=1.038164376506281; "L"=1.03876375806075; "M"=1.034859309580701; "N"=1.038746836432103 }
    8 = @{ "B"=1.02; "C"=1.031421489914182; "D"=1.034709742729889; "E"=1.03504207465137; "F"=1.030230552152257; "I"=1.036430224951466; "J"=1.036438124476844; "K"=1.037444870703894; "L"=1.037776271767409; "M"=1.032978289868451; "N"=1.037909985494715 }
    9 = @{ "B"=1.02; "C"=1.029116188701491; "D"=1.03299713005511; "E"=1.032859098882809; "F"=1.0264633371033; "I"=1.035773413349372; "J"=1.034958941396355; "K"=1.036169534798001; "L"=1.03603195684209; "M"=1.029657341800363; "N"=1.036428701804587 }
    10 = @{ "B"=1.02; "C"=1.027577788727478; "D"=1.031853785157397; "E"=1.031404095381022; "F"=1.023949545485956; "I"=1.035327630835807; "J"=1.033968319311888; "K"=1.035314263372312; "L"=1.034866195349131; "M"=1.02743875873737; "N"=1.035436672923135 }
    11 = @{ "B"=1.02; "C"=1.026911267208977; "D"=1.031358321818441; "E"=1.030774129266148; "F"=1.022860407167323; "I"=1.035132722549318; "J"=1.033538292367316; "K"=1.034942718574033; "L"=1.03436071472696; "M"=1.026476914783653; "N"=1.035006035290977 }
    12 = @{ "B"=1.02; "C"=1.026663632418259; "D"=1.031174226171839; "E"=1.030540140024408; "F"=1.022455747263852; "I"=1.03506004185223; "J"=1.033378397904916; "K"=1.03480452843364; "L"=1.034172850964928; "M"=1.026119458638381; "N"=1.034845913760099 }
    13 = @{ "B"=1.02; "C"=1.026716753613386; "D"=1.031213717965103; "E"=1.030590331151811; "F"=1.022542553083564; "I"=1.035075644925045; "J"=1.033412703214364; "K"=1.034834178926729; "L"=1.034213153178966; "M"=1.02619614265181; "N"=1.034880267787022 }
    14 = @{ "B"=1.02; "C"=1.026890798851531; "D"=1.031343105610015; "E"=1.030754787479319; "F"=1.022826960049566; "I"=1.035126720517353; "J"=1.033525078777816; "K"=1.034931299447166; "L"=1.03434518801209; "M"=1.02644737114966; "N"=1.034992802936664 }
    15 = @{ "B"=1.02; "C"=1.026998026006834; "D"=1.031422817847376; "E"=1.030856115516115; "F"=1.023002178362701; "I"=1.035158152352288; "J"=1.033594295404502; "K"=1.034991114489648; "L"=1.034426524991487; "M"=1.026602136671115; "N"=1.035062117858899 }
    16 = @{ "B"=1.02; "C"=1.027622015360307; "D"=1.031886659239615; "E"=1.031445905354159; "F"=1.024021813651321; "I"=1.035340526581052; "J"=1.033996835931522; "K"=1.03533889610561; "L"=1.03489972765016; "M"=1.027502567654898; "N"=1.03546523003964 }
    17 = @{ "B"=1.02; "C"=1.028013323624294; "D"=1.032177510756262; "E"=1.031815880446787; "F"=1.024661224432206; "I"=1.035454421097223; "J"=1.034249048915382; "K"=1.035556726872336; "L"=1.035196367647184; "M"=1.028067062776425; "N"=1.035717801194866 }
    18 = @{ "B"=1.02; "C"=1.028241530134688; "D"=1.032347122161675; "E"=1.032031686502042; "F"=1.025034119735508; "I"=1.035520672329922; "J"=1.034396056237067; "K"=1.035683667490362; "L"=1.035369325430198; "M"=1.028396209558266; "N"=1.03586501728381 }
    19 = @{ "B"=1.02; "C"=1.02831933642446; "D"=1.032404948938975; "E"=1.032105271804983; "F"=1.025161257085796; "I"=1.035543231512899; "J"=1.034446164280076; "K"=1.03572693123755; "L"=1.035428288173747; "M"=1.028508421029522; "N"=1.035915196485987 }
    20 = @{ "B"=1.02; "C"=1.027971343785014; "D"=1.032146309007731; "E"=1.031776185021876; "F"=1.02459262819447; "I"=1.03544222007094; "J"=1.034221999637915; "K"=1.035533367734624; "L"=1.035164547950888; "M"=1.028006509597962; "N"=1.035690713504322 }
    21 = @{ "B"=1.02; "C"=1.026839548506671; "D"=1.03130500579083; "E"=1.03070635896927; "F"=1.022743212260087; "I"=1.035111687846144; "J"=1.033491991490596; "K"=1.034902704899448; "L"=1.034306309970822; "M"=1.026373395765654; "N"=1.034959668661701 }
    22 = @{ "B"=1.02; "C"=1.026127601759306; "D"=1.030775706439928; "E"=1.03003376404136; "F"=1.021579797142; "I"=1.034902231093274; "J"=1.033032061094005; "K"=1.034505129670079; "L"=1.033766089407381; "M"=1.025345521592694; "N"=1.034499085111192 }
    23 = @{ "B"=1.02; "C"=1.026505051069101; "D"=1.031056330238916; "E"=1.030390315134848; "F"=1.022196606312865; "I"=1.035013423483279; "J"=1.033275968803982; "K"=1.034715991779107; "L"=1.034052528875126; "M"=1.025890520897437; "N"=1.034743339198094 }
    24 = @{ "B"=1.02; "C"=1.027990312781238; "D"=1.032160407850079; "E"=1.031794121657365; "F"=1.024623624071962; "I"=1.035447733750207; "J"=1.034234222363877; "K"=1.035543923082845; "L"=1.035178926109649; "M"=1.028033871322469; "N"=1.035702953587957 }
    25 = @{ "B"=1.02; "C"=1.029712430149042; "D"=1.033440163624249; "E"=1.033423392835367; "F"=1.027437630959822; "I"=1.035944607643391; "J"=1.035342135994919; "K"=1.036500127717948; "L"=1.036483409899328; "M"=1.030516674805242; "N"=1.036812440583434 }
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}